# Adds two more "Docente(s) Responsável(eis)" entries after the existing
# "101761 - Arnaldo Márcio Ramalho Prata" line, each on its own line
# (separated by manual line breaks <w:br/>, matching the pattern already
# used elsewhere in this document), by rebuilding that single list
# paragraph's OOXML in place.

$d = $word.ActiveDocument

# Locate the paragraph that currently holds only the first professor's name.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "101761 - Arnaldo Márcio Ramalho Prata") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>101761 - Arnaldo Márcio Ramalho Prata</w:t><w:br/></w:r><w:r><w:t>6007846 - Júlio César dos Santos</w:t><w:br/></w:r><w:r><w:t>1814052 - Silvio Silverio da Silva</w:t></w:r></w:p>
'@

$target.Range.InsertXML($newParagraphXml) | Out-Null
